$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.063.45"
$ws.Range("E2").Value = "  +0.20%  "

# Row 3
$ws.Range("D3").Value = "3.309.08"
$ws.Range("E3").Value = "  -0.49%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'187.55"
$ws.Range("E5").Value = "  +3.73%  "

# Row 6
$ws.Range("D6").Value = "'556.95"
$ws.Range("E6").Value = "  -0.24%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "  -0.20%  "

# Row 9
$ws.Range("D9").Value = "3.304.76"
$ws.Range("E9").Value = "  -0.46%  "

# Row 10
$ws.Range("E10").Value = "  -0.50%  "

# Row 11
$ws.Range("E11").Value = "  +0.74%  "

# Row 12
$ws.Range("D12").Value = "'47.42"
$ws.Range("E12").Value = "  +0.69%  "

# Row 13
$ws.Range("E13").Value = "  +2.83%  "

# Row 14
$ws.Range("D14").Value = "'8.71"
$ws.Range("E14").Value = "  +2.35%  "

# Row 15
$ws.Range("D15").Value = "3.847.53"
$ws.Range("E15").Value = "  -0.38%  "

# Row 16
$ws.Range("D16").Value = "'604.09"

# Row 17
$ws.Range("D17").Value = "66.161.87"
$ws.Range("E17").Value = "  +0.21%  "

# Row 18
$ws.Range("D18").Value = "'17.98"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19
$ws.Range("E19").Value = "  +1.15%  "

# Row 20
$ws.Range("D20").Value = "3.321.21"
$ws.Range("E20").Value = "  +0.08%  "

# Row 21
$ws.Range("D21").Value = "'11.10"
$ws.Range("E21").Value = "  -2.21%  "

# Row 22
$ws.Range("D22").Value = "'0.910"
$ws.Range("E22").Value = "  +1.05%  "

# Row 23
$ws.Range("D23").Value = "'18.47"
$ws.Range("E23").Value = "  +10.40%  "

# Row 24
$ws.Range("D24").Value = "'5.12"
$ws.Range("E24").Value = "  +1.76%  "

# Row 25
$ws.Range("D25").Value = "'100.06"
$ws.Range("E25").Value = "  +0.36%  "

# Row 26
$ws.Range("D26").Value = "'3.96"
$ws.Range("E26").Value = "  -0.68%  "

# Row 27
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'2.77"
$ws.Range("E27").Value = "  +5.16%  "

# Row 28
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'5.98"
$ws.Range("E28").Value = "  -0.20%  "

# Row 29
$ws.Range("D29").Value = "'9.61"
$ws.Range("E29").Value = "  +4.29%  "

# Row 30
$ws.Range("D30").Value = "'8.63"
$ws.Range("E30").Value = "  +0.14%  "

# Row 31
$ws.Range("D31").Value = "'30.35"
$ws.Range("E31").Value = "  -0.47%  "

# Row 32
$ws.Range("D32").Value = "'6.75"
$ws.Range("E32").Value = "  +8.99%  "

# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'583.48"
$ws.Range("E33").Value = "  +10.15%  "

# Row 34
$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "'3.79"
$ws.Range("E34").Value = "  +0.88%  "

# Row 35
$ws.Range("D35").Value = "'11.13"
$ws.Range("E35").Value = "  +1.53%  "

# Row 36
$ws.Range("E36").Value = "  +1.35%  "

# Row 37
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'57.33"
$ws.Range("E37").Value = "  -0.81%  "

# Row 38
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "3.709.29"
$ws.Range("E38").Value = "  -1.24%  "

# Row 39
$ws.Range("E39").Value = "  -0.08%  "

# Row 40
$ws.Range("D40").Value = "'3.49"
$ws.Range("E40").Value = "  +16.12%  "

# Row 41
$ws.Range("E41").Value = "  +5.84%  "

# Row 42
$ws.Range("E42").Value = "  +2.17%  "

# Row 43
$ws.Range("D43").Value = "'33.91"
$ws.Range("E43").Value = "  +6.95%  "

# Row 44
$ws.Range("E44").Value = "  -4.61%  "

# Row 45
$ws.Range("D45").Value = "'2.67"
$ws.Range("E45").Value = "  +1.27%  "

# Row 46
$ws.Range("D46").Value = "'0.342"
$ws.Range("E46").Value = "  +1.15%  "

# Row 47
$ws.Range("E47").Value = "  +3.21%  "

# Row 48
$ws.Range("D48").Value = "'0.0422"
$ws.Range("E48").Value = "  +2.65%  "

# Row 49
$ws.Range("E49").Value = "  +0.63%  "

# Row 50
$ws.Range("D50").Value = "'2.59"
$ws.Range("E50").Value = "  -0.02%  "

# Row 51
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.04%  "
